$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data changes: update Counts (cps) and Error (cps) for RaFHYASW_2A, RaFHYASW_2B, RaFHYASW_2C
$ws.Range("B14").Value = 4.942222222222222
$ws.Range("C14").Value = 0.1815442962962963

$ws.Range("B15").Value = 4.836111111111111
$ws.Range("C15").Value = 0.179629287037037

$ws.Range("B16").Value = 4.882777777777778
$ws.Range("C16").Value = 0.1804186388888889
